$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1185.5
$ws.Range("I19").Value = 499.6
$ws.Range("J19").Value = 1675.4286
$ws.Range("K19").Value = 499.6
$ws.Range("L19").Value = 1675.4286
$ws.Range("M19").Value = -324.6
$ws.Range("N19").Value = -2025.4286

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1639.0834
$ws.Range("I41").Value = 1961.6666
$ws.Range("K41").Value = 1961.6666
$ws.Range("M41").Value = -1521.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2490.4707
$ws.Range("I137").Value = 2006.7858
$ws.Range("J137").Value = 2829.05
$ws.Range("K137").Value = 6020.357400000001
$ws.Range("L137").Value = 8487.150000000001
$ws.Range("M137").Value = -3470.357400000001
$ws.Range("N137").Value = -13587.15

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3863.8
$ws.Range("J138").Value = 2668
$ws.Range("L138").Value = 8004
$ws.Range("N138").Value = -18284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5162.378
$ws.Range("I32").Value = 3968.027
$ws.Range("K32").Value = 3968.027
$ws.Range("M32").Value = -3681.027

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1511.5555
$ws.Range("I61").Value = 1454.2667
$ws.Range("K61").Value = 1454.2667
$ws.Range("M61").Value = -1242.2667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1031.275
$ws.Range("I74").Value = 538.96875
$ws.Range("J74").Value = 3000.5
$ws.Range("K74").Value = 538.96875
$ws.Range("L74").Value = 3000.5
$ws.Range("M74").Value = 335.03125
$ws.Range("N74").Value = -4748.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1031.275
$ws.Range("I77").Value = 538.96875
$ws.Range("J77").Value = 3000.5
$ws.Range("K77").Value = 2694.84375
$ws.Range("L77").Value = 15002.5
$ws.Range("M77").Value = 1673.15625
$ws.Range("N77").Value = -23738.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1511.5555
$ws.Range("I136").Value = 1454.2667
$ws.Range("K136").Value = 4362.800099999999
$ws.Range("M136").Value = -1812.800099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2802.6667
$ws.Range("I31").Value = 2663.889
$ws.Range("K31").Value = 2663.889
$ws.Range("M31").Value = -2368.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2802.6667
$ws.Range("I34").Value = 2663.889
$ws.Range("K34").Value = 2663.889
$ws.Range("M34").Value = -2461.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3722
$ws.Range("I86").Value = 3333
$ws.Range("K86").Value = 3333
$ws.Range("M86").Value = -2210

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3722
$ws.Range("I89").Value = 3333
$ws.Range("K89").Value = 16665
$ws.Range("M89").Value = -11049

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3216.75
$ws.Range("I99").Value = 2146.8
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 2146.8
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -648.8000000000002
$ws.Range("N99").Value = -7996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1090.5
$ws.Range("I105").Value = 1103.4286
$ws.Range("K105").Value = 1103.4286
$ws.Range("M105").Value = 643.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3216.75
$ws.Range("I126").Value = 2146.8
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 6440.400000000001
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -3970.400000000001
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3903.375
$ws.Range("I132").Value = 3052.2222
$ws.Range("J132").Value = 4997.7144
$ws.Range("K132").Value = 9156.6666
$ws.Range("L132").Value = 14993.1432
$ws.Range("M132").Value = -6626.6666
$ws.Range("N132").Value = -20053.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3230.611
$ws.Range("I134").Value = 2781.1538
$ws.Range("K134").Value = 8343.4614
$ws.Range("M134").Value = -5808.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 151.44444
$ws.Range("I23").Value = 241.25
$ws.Range("K23").Value = 723.75
$ws.Range("M23").Value = -488.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 70000
$ws.Range("J37").Value = 70000
$ws.Range("L37").Value = 210000
$ws.Range("N37").Value = -210224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 21889.666
$ws.Range("J75").Value = 26213.6
$ws.Range("L75").Value = 78640.79999999999
$ws.Range("N75").Value = -80636.79999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 21889.666
$ws.Range("J78").Value = 26213.6
$ws.Range("L78").Value = 235922.4
$ws.Range("N78").Value = -245906.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 502
$ws.Range("J97").Value = 502
$ws.Range("L97").Value = 1506
$ws.Range("N97").Value = -2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 709.7
$ws.Range("J117").Value = 669.1667
$ws.Range("L117").Value = 2007.5001
$ws.Range("N117").Value = -8891.500099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19093.79
$ws.Range("J131").Value = 20665.258
$ws.Range("L131").Value = 61995.774
$ws.Range("N131").Value = -72075.774

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2235.3215
$ws.Range("I140").Value = 1367.3334
$ws.Range("J140").Value = 3236.8462
$ws.Range("K140").Value = 4102.0002
$ws.Range("L140").Value = 9710.5386
$ws.Range("M140").Value = 1077.9998
$ws.Range("N140").Value = -20070.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4339.5
$ws.Range("I102").Value = 5204.857
$ws.Range("J102").Value = 3128
$ws.Range("K102").Value = 5204.857
$ws.Range("L102").Value = 3128
$ws.Range("M102").Value = -3582.857
$ws.Range("N102").Value = -6372

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1258.381
$ws.Range("I122").Value = 1287.2632
$ws.Range("K122").Value = 3861.7896
$ws.Range("M122").Value = -1411.7896

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3241.1072
$ws.Range("I132").Value = 2669.7778
$ws.Range("K132").Value = 8009.3334
$ws.Range("M132").Value = -5479.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 50395.668
$ws.Range("J138").Value = 50395.668
$ws.Range("L138").Value = 50395.668
$ws.Range("N138").Value = -60675.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = 33516.5
$ws.Range("J38").Value = 33516.5
$ws.Range("L38").Value = 33516.5
$ws.Range("N38").Value = -34336.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2360.8
$ws.Range("I68").Value = 1701
$ws.Range("K68").Value = 1701
$ws.Range("M68").Value = -952

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2360.8
$ws.Range("I71").Value = 1701
$ws.Range("K71").Value = 8505
$ws.Range("M71").Value = -4761

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2557
$ws.Range("I100").Value = 1699.2
$ws.Range("K100").Value = 1699.2
$ws.Range("M100").Value = -1158.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7759.2856
$ws.Range("I132").Value = 10299.333
$ws.Range("J132").Value = 5854.25
$ws.Range("K132").Value = 30897.999
$ws.Range("L132").Value = 17562.75
$ws.Range("M132").Value = -28367.999
$ws.Range("N132").Value = -22622.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4788.609
$ws.Range("I136").Value = 3750.8
$ws.Range("J136").Value = 6734.5
$ws.Range("K136").Value = 11252.4
$ws.Range("L136").Value = 20203.5
$ws.Range("M136").Value = -8702.400000000001
$ws.Range("N136").Value = -25303.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1197.8889
$ws.Range("I81").Value = 641
$ws.Range("J81").Value = 1894
$ws.Range("K81").Value = 1282
$ws.Range("L81").Value = 3788
$ws.Range("M81").Value = -221
$ws.Range("N81").Value = -5910

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1197.8889
$ws.Range("I84").Value = 641
$ws.Range("J84").Value = 1894
$ws.Range("K84").Value = 6410
$ws.Range("L84").Value = 18940
$ws.Range("M84").Value = -1106
$ws.Range("N84").Value = -29548

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 633.4545000000001
$ws.Range("I107").Value = 296.05884
$ws.Range("K107").Value = 888.17652
$ws.Range("M107").Value = 1031.82348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 112489.57
$ws.Range("I122").Value = 195456.75
$ws.Range("K122").Value = 586370.25
$ws.Range("M122").Value = -583920.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1524.9231
$ws.Range("I132").Value = 864.1667
$ws.Range("K132").Value = 2592.5001
$ws.Range("M132").Value = -62.5001000000002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1746.5555
$ws.Range("I136").Value = 1969.1
$ws.Range("K136").Value = 5907.299999999999
$ws.Range("M136").Value = -3357.299999999999
